$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.713.92'
$ws.Range('E2').Value = '  +0.30%  '
$ws.Range('D3').Value = '1.599.28'
$ws.Range('E3').Value = '  +0.16%  '
$ws.Range('E4').Value = '  +0.33%  '
$ws.Range('D5').Value = "'211.46"
$ws.Range('E5').Value = '  -0.11%  '
$ws.Range('E6').Value = '  -0.78%  '
$ws.Range('E7').Value = '  +0.33%  '
$ws.Range('E8').Value = '  +0.07%  '
$ws.Range('E9').Value = '  +0.66%  '
$ws.Range('D10').Value = "'19.54"
$ws.Range('E10').Value = '  +0.10%  '
$ws.Range('D11').Value = "'0.0842"
$ws.Range('E11').Value = '  +0.46%  '
$ws.Range('D12').Value = '1.824.50'
$ws.Range('E12').Value = '  +0.22%  '
$ws.Range('D13').Value = '1.592.75'
$ws.Range('E13').Value = '  +0.16%  '
$ws.Range('E14').Value = '  +0.55%  '
$ws.Range('E15').Value = '  +0.12%  '
$ws.Range('D16').Value = "'65.29"
$ws.Range('E16').Value = '  +1.22%  '
$ws.Range('D17').Value = '26.686.00'
$ws.Range('E17').Value = '  +0.26%  '
$ws.Range('D18').Value = '0.0₃0752'
$ws.Range('E18').Value = '  +2.80%  '
$ws.Range('D19').Value = "'7.25"
$ws.Range('E19').Value = '  +4.07%  '
$ws.Range('E20').Value = '  +0.34%  '
$ws.Range('D21').Value = "'209.14"
$ws.Range('E21').Value = '  +0.25%  '
$ws.Range('E22').Value = '  +0.29%  '
$ws.Range('E23').Value = '  +0.06%  '
$ws.Range('E24').Value = '  +0.48%  '
$ws.Range('D25').Value = "'142.36"
$ws.Range('E25').Value = '  -1.99%  '
$ws.Range('E26').Value = '  +0.38%  '
$ws.Range('E27').Value = '  -0.38%  '
$ws.Range('D28').Value = "'0.115"
$ws.Range('E28').Value = '  +0.11%  '
$ws.Range('D29').Value = "'15.35"
$ws.Range('E29').Value = '  +0.56%  '
$ws.Range('E30').Value = '  +2.96%  '
$ws.Range('E31').Value = '  -0.32%  '
$ws.Range('E32').Value = '  +0.65%  '
$ws.Range('D33').Value = "'2.97"
$ws.Range('E33').Value = '  +1.50%  '
$ws.Range('D34').Value = '1.292.26'
$ws.Range('E34').Value = '  +0.88%  '
$ws.Range('D35').Value = "'0.622"
$ws.Range('E35').Value = '  -4.86%  '
$ws.Range('E36').Value = '  +1.09%  '
$ws.Range('E37').Value = '  -0.09%  '
$ws.Range('E38').Value = '  -0.09%  '
$ws.Range('E39').Value = '  +20.89%  '
$ws.Range('E40').Value = '  -2.09%  '
$ws.Range('E41').Value = '  -0.86%  '
$ws.Range('D42').Value = "'2.20"
$ws.Range('E42').Value = '  +0.18%  '
$ws.Range('E43').Value = '  -0.23%  '
$ws.Range('D44').Value = "'63.15"
$ws.Range('E44').Value = '  -2.05%  '
$ws.Range('D45').Value = '1.736.75'
$ws.Range('E45').Value = '  +0.19%  '
$ws.Range('D46').Value = "'91.24"
$ws.Range('E46').Value = '  +1.67%  '
$ws.Range('E47').Value = '  -1.69%  '
$ws.Range('E48').Value = '  -1.85%  '
$ws.Range('E49').Value = '  +0.50%  '
$ws.Range('E50').Value = '  +0.19%  '
$ws.Range('D51').Value = "'7.38"
$ws.Range('E51').Value = '  -1.18%  '
